$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before "name" (column C) for prolificid
$ws.Columns.Item(3).Insert()

# Header row
$ws.Range("C1").Value = "prolificid"

# Data rows: set prolificid (C), name (D), realeffort (F), race (G) and fix index (B) order/values
$ws.Range("B2").Value = 44
$ws.Range("C2").Value = "60b091ed11ccda59e3fc7761"
$ws.Range("D2").Value = "Myles"
$ws.Range("F2").Value = 13.48524091344236
$ws.Range("G2").Value = "Black or African American"

$ws.Range("B3").Value = 3
$ws.Range("C3").Value = "601d69a993d94008fb2b25dc"
$ws.Range("D3").Value = "Quinterius"
$ws.Range("F3").Value = 8.41655457137187
$ws.Range("G3").Value = "Black or African American"

$ws.Range("B4").Value = 30
$ws.Range("C4").Value = "60c2341fe95d71ee52c043f0"
$ws.Range("D4").Value = "Matthew"
$ws.Range("F4").Value = 7.329976078078919
$ws.Range("G4").Value = "White"

$ws.Range("B5").Value = 27
$ws.Range("C5").Value = "5ff8ad350d084e10f500e48a"
$ws.Range("D5").Value = "Drew"
$ws.Range("F5").Value = 7.217939142941571
$ws.Range("G5").Value = "White"

$ws.Range("B6").Value = 22
$ws.Range("C6").Value = "60db4fde6193c50664c9c478"
$ws.Range("D6").Value = "Edosagbe"
$ws.Range("F6").Value = 5.366656653952606
$ws.Range("G6").Value = "Black or African American"

$ws.Range("B7").Value = 26
$ws.Range("C7").Value = "5dd671942b033b5ec8bc97b4"
$ws.Range("D7").Value = "Juan"
$ws.Range("F7").Value = 5.23281095267766
$ws.Range("G7").Value = "Hispanic"

$ws.Range("B8").Value = 32
$ws.Range("C8").Value = "60bf9943e4e04642d4634ecc"
$ws.Range("D8").Value = "Jamarii"
$ws.Range("F8").Value = 5.099160166839549
$ws.Range("G8").Value = "Black or African American"

$ws.Range("B9").Value = 33
$ws.Range("C9").Value = "60b322994d0b901954690036"
$ws.Range("D9").Value = "Brennan"
$ws.Range("F9").Value = 4.111111164481627
$ws.Range("G9").Value = "White"

$ws.Range("B10").Value = 2
$ws.Range("C10").Value = "5e2522d6b734b47915f88275"
$ws.Range("D10").Value = "Corey"
$ws.Range("F10").Value = 4.088873157346726
$ws.Range("G10").Value = "White"

$ws.Range("B11").Value = 49
$ws.Range("C11").Value = "6088fc724afd5c008db33e9d"
$ws.Range("D11").Value = "Masuf"
$ws.Range("F11").Value = 3.326168639869025
$ws.Range("G11").Value = "Asian"

$ws.Range("B12").Value = 50
$ws.Range("C12").Value = "6097b95056caf5ebb2720002"
$ws.Range("D12").Value = "Damian"
$ws.Range("F12").Value = 2.338825072770427
$ws.Range("G12").Value = "Black or African American"

$ws.Range("B13").Value = 29
$ws.Range("C13").Value = "60b83826821417f8e484a207"
$ws.Range("D13").Value = "Eli"
$ws.Range("F13").Value = 2.32406207211523
$ws.Range("G13").Value = "White"
